$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Simple single-value cell replacements near the top of the table ---
$t.Rows(1).Cells(1).Range.Text  = "0M"
$t.Rows(2).Cells(1).Range.Text  = "0M"
$t.Rows(3).Cells(1).Range.Text  = "0M"
$t.Rows(4).Cells(1).Range.Text  = "804"
$t.Rows(5).Cells(1).Range.Text  = "0.00002"
$t.Rows(6).Cells(1).Range.Text  = "0.00008"

# --- Remove the extra "0.00004" row (the third of three consecutive
#     "0.00004" rows, immediately before the "0.00005" row) ---
$t.Rows(10).Delete()

# After the delete, the former row 12 ("0.00435") is now row 11.
$t.Rows(11).Cells(1).Range.Text = "0.00007"

# Insert a brand-new row right after it (i.e. before the now-row-12
# "100.0" row) holding the new value "0.03735".
$newRow = $t.Rows.Add($t.Rows(12))
$newRow.Cells(1).Range.Text = "0.03735"

# --- Collapse the three multi-value (tab-separated) rows near the
#     bottom of the table down to the single leftover summary value ---
$t.Rows(44).Cells(1).Range.Text = "99.99"
$t.Rows(45).Cells(1).Range.Text = "0.04"
$t.Rows(46).Cells(1).Range.Text = "545"
